# "Chiffres COVID-19 Valais" - daily data upload.
# Updates the raw input columns (C, E, F, G, L, M) for a handful of rows
# near the end of the data table (rows 271-279, dates late Nov/early Dec
# 2020). Columns B, H, J, K are live volatile formulas (cumulative
# totals / sums depending on TODAY()) and are left untouched - they
# recompute automatically from the inputs below once Excel recalculates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumericValue {
    # Columns L ("Nb nouveaux décès à l'hôpital") and M ("Nb nouveaux
    # décès extra-hospitaliers") are formatted as Text (numFmt "@") in
    # this sheet. Writing a plain number straight into .Value on a
    # Text-formatted cell makes Excel store it as a text string instead
    # of a number (matches real Excel behaviour for Text-formatted
    # cells). The existing data in those columns is genuinely numeric,
    # so flip the format to a plain numeric one, write the value, then
    # restore the original Text format - this keeps both the stored
    # type (number) and the visible style (s="17"/"18", numFmt "@")
    # identical to the rest of the column.
    param($cell, $value)
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "0"
    $cell.Value = $value
    $cell.NumberFormat = $fmt
}

# Row 271 (2020-11-23): one fewer intubated, one more hospitalised hors SI.
$ws.Cells.Item(271, 6).Value = 24    # F271 Patients intubés
$ws.Cells.Item(271, 7).Value = 189   # G271 Hospitalisés hors SI

# Row 272 (2020-11-24)
$ws.Cells.Item(272, 6).Value = 26    # F272
$ws.Cells.Item(272, 7).Value = 196   # G272

# Row 273 (2020-11-25)
$ws.Cells.Item(273, 6).Value = 25    # F273
$ws.Cells.Item(273, 7).Value = 191   # G273

# Row 274 (2020-11-26)
$ws.Cells.Item(274, 3).Value = 145   # C274 Nb nouveaux cas positifs
$ws.Cells.Item(274, 6).Value = 24    # F274
$ws.Cells.Item(274, 7).Value = 188   # G274
Set-NumericValue $ws.Cells.Item(274, 12) 4   # L274 Nb nouveaux décès à l'hôpital

# Row 275 (2020-11-27)
$ws.Cells.Item(275, 3).Value = 143   # C275
$ws.Cells.Item(275, 5).Value = 28    # E275 Patients SI total
$ws.Cells.Item(275, 6).Value = 20    # F275
$ws.Cells.Item(275, 7).Value = 173   # G275
Set-NumericValue $ws.Cells.Item(275, 12) 4   # L275
Set-NumericValue $ws.Cells.Item(275, 13) 1   # M275 Nb nouveaux décès extra-hospitaliers

# Row 276 (2020-11-28)
$ws.Cells.Item(276, 3).Value = 159   # C276
$ws.Cells.Item(276, 5).Value = 28    # E276
$ws.Cells.Item(276, 6).Value = 22    # F276
$ws.Cells.Item(276, 7).Value = 165   # G276
Set-NumericValue $ws.Cells.Item(276, 12) 3   # L276
Set-NumericValue $ws.Cells.Item(276, 13) 1   # M276

# Row 277 (2020-11-29) - previously blank, now filled in with real figures.
$ws.Cells.Item(277, 3).Value = 77    # C277
$ws.Cells.Item(277, 5).Value = 28    # E277
$ws.Cells.Item(277, 6).Value = 21    # F277
$ws.Cells.Item(277, 7).Value = 155   # G277
Set-NumericValue $ws.Cells.Item(277, 12) 0   # L277
Set-NumericValue $ws.Cells.Item(277, 13) 0   # M277

# Row 278 (2020-11-30) - previously blank, now filled in.
$ws.Cells.Item(278, 3).Value = 44    # C278
$ws.Cells.Item(278, 5).Value = 27    # E278
$ws.Cells.Item(278, 6).Value = 20    # F278
$ws.Cells.Item(278, 7).Value = 152   # G278
Set-NumericValue $ws.Cells.Item(278, 12) 3   # L278
Set-NumericValue $ws.Cells.Item(278, 13) 0   # M278

# Row 279 (2020-12-01) - previously blank, now filled in.
$ws.Cells.Item(279, 3).Value = 26    # C279
$ws.Cells.Item(279, 5).Value = 27    # E279
$ws.Cells.Item(279, 6).Value = 21    # F279
$ws.Cells.Item(279, 7).Value = 159   # G279
Set-NumericValue $ws.Cells.Item(279, 12) 0   # L279
Set-NumericValue $ws.Cells.Item(279, 13) 0   # M279

# Row 280 stays without raw inputs; its formula cells (B/H/J/K) simply
# pick up the new totals once recalculated - no direct write needed.
